$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: Manager for project "Acacia Breeze" (row 2) should be "Michael", not "T8765432F"
$ws.Range("L2").Value = "Michael"

# Remove the duplicate/erroneous last row (row 6 - "Archipelago")
$ws.Rows(6).Delete()
